$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current last row (old row 34), shifting it down to row 36
$ws.Rows("34:35").Insert()
$ws.Rows("34:35").RowHeight = 15

# Row 34 — new "dataSetup" entry with caseChildLocationSetup script (Consolas font)
$ws.Range("A34").Value = "dataSetup"
$ws.Range("B34").Value = "caseChildLocationSetup"
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = "YES"

# Apply the Consolas font (10pt, black) to the new script-name cell.
# Ordering chosen so each intermediate step reuses an already-existing font entry.
$ws.Range("B34").Font.Color = "#000000"
$ws.Range("B34").Font.Name = "Consolas"
$ws.Range("B34").Font.Size = 10

# Row 35 — new "testT4249" entry
$ws.Range("A35").Value = 20
$ws.Range("B35").Value = "testT4249"
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = "YES"

# Row 36 (previously row 34) — script name updated from testT4245 to testT4145
$ws.Range("B36").Value = "testT4145"

# Update the view: scroll down one row and move the active selection to G34
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G34").Select()
